# Natmi following Dr Hou advice
# Update ligand/receptor-expressing cell counts (E, K: 1 -> 3) and
# recompute the dependent expression/specificity metrics for rows 2-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> @{ Column = NewValue } map, one entry per data row (2..17)
$rowUpdates = @{
    2 = @{ "E" = 3; "G" = 19.95578266666667; "H" = 59.867348; "I" = 0.0117373419656925; "J" = 0.0117373419656925; "K" = 3; "M" = 2.237200333333333; "N" = 6.711601; "O" = 0.1121050933480713; "P" = 0.1121050933480713; "Q" = 44.64508363379423; "R" = 401.805752704148; "S" = 0.001315815816722192; "T" = 0.001315815816722192 }
    3 = @{ "E" = 3; "G" = 19.95578266666667; "H" = 59.867348; "I" = 0.0117373419656925; "J" = 0.0117373419656925; "K" = 3; "M" = 3.153682; "N" = 9.461046; "O" = 0.158029573718759; "P" = 0.158029573718759; "Q" = 62.93419259177867; "R" = 566.4077333260079; "S" = 0.001854847147429686; "T" = 0.001854847147429686 }
    4 = @{ "E" = 3; "G" = 19.95578266666667; "H" = 59.867348; "I" = 0.0117373419656925; "J" = 0.0117373419656925; "K" = 3; "M" = 13.44189533333333; "N" = 40.325686; "O" = 0.6735672745377762; "P" = 0.6735672745377762; "Q" = 268.2435419000809; "R" = 2414.191877100728; "S" = 0.00790588943814936; "T" = 0.00790588943814936 }
    5 = @{ "E" = 3; "G" = 19.95578266666667; "H" = 59.867348; "I" = 0.0117373419656925; "J" = 0.0117373419656925; "K" = 3; "M" = 1.123499666666667; "N" = 3.370499; "O" = 0.05629805839539345; "P" = 0.05629805839539345; "Q" = 22.42031517407245; "R" = 201.782836566652; "S" = 0.0006607895633912583; "T" = 0.0006607895633912583 }
    6 = @{ "E" = 3; "G" = 1637.343343333333; "H" = 4912.03003; "I" = 0.9630320723052701; "J" = 0.9630320723052702; "K" = 3; "M" = 2.237200333333333; "N" = 6.711601; "O" = 0.1121050933480713; "P" = 0.1121050933480713; "Q" = 3663.065073486448; "R" = 32967.58566137803; "S" = 0.1079608003629689; "T" = 0.1079608003629689 }
    7 = @{ "E" = 3; "G" = 1637.343343333333; "H" = 4912.03003; "I" = 0.9630320723052701; "J" = 0.9630320723052702; "K" = 3; "M" = 3.153682; "N" = 9.461046; "O" = 0.158029573718759; "P" = 0.158029573718759; "Q" = 5163.660229690153; "R" = 46472.94206721138; "S" = 0.1521875478638949; "T" = 0.152187547863895 }
    8 = @{ "E" = 3; "G" = 1637.343343333333; "H" = 4912.03003; "I" = 0.9630320723052701; "J" = 0.9630320723052702; "K" = 3; "M" = 13.44189533333333; "N" = 40.325686; "O" = 0.6735672745377762; "P" = 0.6735672745377762; "Q" = 22008.99784581673; "R" = 198080.9806123506; "S" = 0.6486668882351274; "T" = 0.6486668882351275 }
    9 = @{ "E" = 3; "G" = 1637.343343333333; "H" = 4912.03003; "I" = 0.9630320723052701; "J" = 0.9630320723052702; "K" = 3; "M" = 1.123499666666667; "N" = 3.370499; "O" = 0.05629805839539345; "P" = 0.05629805839539345; "Q" = 1839.554700453886; "R" = 16555.99230408497; "S" = 0.05421683584327886; "T" = 0.05421683584327887 }
    10 = @{ "E" = 3; "G" = 17.50081933333334; "H" = 52.502458; "I" = 0.01029341242216722; "J" = 0.01029341242216722; "K" = 3; "M" = 2.237200333333333; "N" = 6.711601; "O" = 0.1121050933480713; "P" = 0.1121050933480713; "Q" = 39.15283884613979; "R" = 352.375549615258; "S" = 0.001153943960457253; "T" = 0.001153943960457253 }
    11 = @{ "E" = 3; "G" = 17.50081933333334; "H" = 52.502458; "I" = 0.01029341242216722; "J" = 0.01029341242216722; "K" = 3; "M" = 3.153682; "N" = 9.461046; "O" = 0.158029573718759; "P" = 0.158029573718759; "Q" = 55.19201891678534; "R" = 496.728170251068; "S" = 0.001626663577186464; "T" = 0.001626663577186465 }
    12 = @{ "E" = 3; "G" = 17.50081933333334; "H" = 52.502458; "I" = 0.01029341242216722; "J" = 0.01029341242216722; "K" = 3; "M" = 13.44189533333333; "N" = 40.325686; "O" = 0.6735672745377762; "P" = 0.6735672745377762; "Q" = 235.2441817262432; "R" = 2117.197635536188; "S" = 0.006933305750892463; "T" = 0.006933305750892464 }
    13 = @{ "E" = 3; "G" = 17.50081933333334; "H" = 52.502458; "I" = 0.01029341242216722; "J" = 0.01029341242216722; "K" = 3; "M" = 1.123499666666667; "N" = 3.370499; "O" = 0.05629805839539345; "P" = 0.05629805839539345; "Q" = 19.66216468739356; "R" = 176.959482186542; "S" = 0.0005794991336310383; "T" = 0.0005794991336310384 }
    14 = @{ "E" = 3; "G" = 25.39612333333333; "H" = 76.18836999999999; "I" = 0.01493717330687017; "J" = 0.01493717330687017; "K" = 3; "M" = 2.237200333333333; "N" = 6.711601; "O" = 0.1121050933480713; "P" = 0.1121050933480713; "Q" = 56.81621558670778; "R" = 511.3459402803699; "S" = 0.001674533207922999; "T" = 0.001674533207922999 }
    15 = @{ "E" = 3; "G" = 25.39612333333333; "H" = 76.18836999999999; "I" = 0.01493717330687017; "J" = 0.01493717330687017; "K" = 3; "M" = 3.153682; "N" = 9.461046; "O" = 0.158029573718759; "P" = 0.158029573718759; "Q" = 80.09129702611332; "R" = 720.82167323502; "S" = 0.002360515130247918; "T" = 0.002360515130247919 }
    16 = @{ "E" = 3; "G" = 25.39612333333333; "H" = 76.18836999999999; "I" = 0.01493717330687017; "J" = 0.01493717330687017; "K" = 3; "M" = 13.44189533333333; "N" = 40.325686; "O" = 0.6735672745377762; "P" = 0.6735672745377762; "Q" = 341.3720317190911; "R" = 3072.34828547182; "S" = 0.01006119111360696; "T" = 0.01006119111360696 }
    17 = @{ "E" = 3; "G" = 25.39612333333333; "H" = 76.18836999999999; "I" = 0.01493717330687017; "J" = 0.01493717330687017; "K" = 3; "M" = 1.123499666666667; "N" = 3.370499; "O" = 0.05629805839539345; "P" = 0.05629805839539345; "Q" = 28.53253609962556; "R" = 256.79282489663; "S" = 0.000840933855092289; "T" = 0.000840933855092289 }
}

foreach ($rowNum in $rowUpdates.Keys) {
    $colValues = $rowUpdates[$rowNum]
    foreach ($col in $colValues.Keys) {
        $ws.Range("$col$rowNum").Value = $colValues[$col]
    }
}
